# Update "Countries & provincias Spain" workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 22:16"

# Navarra (row 8): Muertes 14 -> 24
$ws.Range("E8").Value = 24

# Rows 19-22: Murcia is inserted (sorted) ahead of Granada, pushing Granada / Zaragoza /
# Gipuzkoa-Guipuzcoa down a row each; Sevilla (row 23) stays where it is.
$ws.Range("A19").Value = "Murcia"
$ws.Range("B19").Value = 345
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 213
$ws.Range("E19").Value = 3

$ws.Range("A20").Value = "Granada"
$ws.Range("B20").Value = 335
$ws.Range("C20").Value = 72
$ws.Range("D20").Value = 276
$ws.Range("E20").Value = 13

$ws.Range("A21").Value = "Zaragoza"
$ws.Range("B21").Value = 329
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 315
$ws.Range("E21").Value = 14

$ws.Range("A22").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B22").Value = 325
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 312
$ws.Range("E22").Value = 13

# Tenerife (row 28): Muertes 5 -> 11
$ws.Range("E28").Value = 11

# Gran Canaria (row 39): Muertes 2 -> 11
$ws.Range("E39").Value = 11

# Fuerteventura (row 56): Muertes 0 -> 11
$ws.Range("E56").Value = 11

# La Palma (row 57): Muertes 0 -> 11
$ws.Range("E57").Value = 11

# Lanzarote (row 58): Muertes 0 -> 11
$ws.Range("E58").Value = 11

# La Gomera (row 61): Muertes 0 -> 11
$ws.Range("E61").Value = 11

# El Hierro (row 62): Muertes 0 -> 11
$ws.Range("E62").Value = 11
